# Munger re-design - BOM - wip
#
# Populate the "Have" (column C) quantities for the parts that are
# already on hand. Column E ("Need") is a shared formula (B-C-D) and
# recalculates automatically once C is filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$haveQty = @{
    3  = 2
    4  = 2
    5  = 2
    6  = 4
    7  = 4
    8  = 1
    9  = 1
    10 = 1
    11 = 4
    12 = 4
    13 = 2
    14 = 4
    15 = 2
    16 = 2
    23 = 2
}

foreach ($row in $haveQty.Keys) {
    $ws.Range("C$row").Value = $haveQty[$row]
}

# Reflect where the user's selection ended up after this editing pass.
$ws.Range("C17").Select()
